# break out stock.yaml completed
# - fix E27 (bsecode) to be a real number instead of text
# - append the newly-scraped row (row 28) to the "10per change" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# E27 was stored as text "544028"; store it as a genuine number.
$ws.Cells.Item(27, 5).Value = 544028

# Append the new scraped row (row 28).
$ws.Cells.Item(28, 1).Value = "21/06/2024 08:44:45"
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(28, 3).Value = "TATATECH"
$ws.Cells.Item(28, 4).Value = "Tata Technologies Ltd"

# bsecode stays text here (matches the still-unconverted "544028" string).
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "544028"

$ws.Cells.Item(28, 6).Value = -0.88
$ws.Cells.Item(28, 7).Value = 1001.45
$ws.Cells.Item(28, 8).Value = 2335783
